# Update cryptocurrency price/volume data (and reorder a few rows)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (avoids Excel auto-numeric coercion
# for values like "1.44" / "0.998" / "0.0772"), then restore the default
# 'Normal' style so no stray number-format style is left on the cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '63.254.26'
Set-TextValue $ws.Range('E2') '  +2.91%  '
Set-TextValue $ws.Range('D3') '3.492.34'
Set-TextValue $ws.Range('E3') '  +2.94%  '
Set-TextValue $ws.Range('E4') '  +0.07%  '
Set-TextValue $ws.Range('D5') '585.15'
Set-TextValue $ws.Range('E5') '  +1.57%  '
Set-TextValue $ws.Range('D6') '148.58'
Set-TextValue $ws.Range('E6') '  +5.57%  '
Set-TextValue $ws.Range('E7') '  -0.12%  '
Set-TextValue $ws.Range('E8') '  +1.08%  '
Set-TextValue $ws.Range('E9') '  +0.20%  '
Set-TextValue $ws.Range('E10') '  +2.73%  '
Set-TextValue $ws.Range('E11') '  +2.62%  '
Set-TextValue $ws.Range('D12') '4.090.13'
Set-TextValue $ws.Range('E12') '  +2.92%  '
Set-TextValue $ws.Range('D13') '29.87'
Set-TextValue $ws.Range('E13') '  +5.75%  '
Set-TextValue $ws.Range('E14') '  -0.21%  '
Set-TextValue $ws.Range('D15') '3.490.48'
Set-TextValue $ws.Range('E15') '  +3.05%  '
Set-TextValue $ws.Range('E16') '  +1.79%  '
Set-TextValue $ws.Range('D17') '63.276.08'
Set-TextValue $ws.Range('E17') '  +2.96%  '
Set-TextValue $ws.Range('D18') '6.31'
Set-TextValue $ws.Range('E18') '  +2.70%  '
Set-TextValue $ws.Range('D19') '14.35'
Set-TextValue $ws.Range('E19') '  +4.99%  '
Set-TextValue $ws.Range('D20') '9.35'
Set-TextValue $ws.Range('E20') '  +4.38%  '
Set-TextValue $ws.Range('D21') '390.26'
Set-TextValue $ws.Range('E21') '  +0.20%  '
Set-TextValue $ws.Range('E22') '  +1.91%  '
Set-TextValue $ws.Range('D23') '75.17'
Set-TextValue $ws.Range('E23') '  -0.27%  '
Set-TextValue $ws.Range('E24') '  +0.00%  '
Set-TextValue $ws.Range('E25') '  +4.21%  '
Set-TextValue $ws.Range('D26') '3.632.32'
Set-TextValue $ws.Range('E26') '  +2.89%  '
Set-TextValue $ws.Range('E27') '  -5.15%  '
Set-TextValue $ws.Range('D28') '7.73'
Set-TextValue $ws.Range('E28') '  +5.91%  '
Set-TextValue $ws.Range('D29') '0.998'
Set-TextValue $ws.Range('E29') '  -0.24%  '
Set-TextValue $ws.Range('D30') '8.29'
Set-TextValue $ws.Range('E30') '  +3.39%  '
Set-TextValue $ws.Range('D31') '1.44'
Set-TextValue $ws.Range('E31') '  +5.26%  '
Set-TextValue $ws.Range('E32') '  +0.10%  '
Set-TextValue $ws.Range('D34') '23.85'
Set-TextValue $ws.Range('E34') '  +1.80%  '
Set-TextValue $ws.Range('D35') '5.37'
Set-TextValue $ws.Range('E35') '  +6.31%  '
Set-TextValue $ws.Range('B36') 'Aptos'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D36') '7.13'
Set-TextValue $ws.Range('E36') '  +2.72%  '
Set-TextValue $ws.Range('B37') 'EnergySwap'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D37') '31.64'
Set-TextValue $ws.Range('E37') '  +21.45%  '
Set-TextValue $ws.Range('B38') 'ImmutableX'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D38') '1.58'
Set-TextValue $ws.Range('E38') '  +6.76%  '
Set-TextValue $ws.Range('B39') 'Monero'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D39') '171.34'
Set-TextValue $ws.Range('E39') '  +2.38%  '
Set-TextValue $ws.Range('D40') '3.527.95'
Set-TextValue $ws.Range('E40') '  +2.97%  '
Set-TextValue $ws.Range('D41') '0.0772'
Set-TextValue $ws.Range('E41') '  +0.34%  '
Set-TextValue $ws.Range('D42') '0.810'
Set-TextValue $ws.Range('E42') '  +4.08%  '
Set-TextValue $ws.Range('B43') 'OKB'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D43') '42.45'
Set-TextValue $ws.Range('E43') '  +0.05%  '
Set-TextValue $ws.Range('B44') 'Filecoin'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D44') '4.49'
Set-TextValue $ws.Range('E44') '  +1.32%  '
Set-TextValue $ws.Range('B45') 'Stacks'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D45') '1.73'
Set-TextValue $ws.Range('E45') '  +3.85%  '
Set-TextValue $ws.Range('D46') '1.21'
Set-TextValue $ws.Range('E46') '  +6.51%  '
Set-TextValue $ws.Range('D47') '2.630.50'
Set-TextValue $ws.Range('E47') '  +7.01%  '
Set-TextValue $ws.Range('B48') 'dogwifhat'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D48') '2.30'
Set-TextValue $ws.Range('E48') '  +10.93%  '
Set-TextValue $ws.Range('B49') 'InjectiveProtocol'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D49') '23.56'
Set-TextValue $ws.Range('E49') '  +2.19%  '
Set-TextValue $ws.Range('D50') '6.79'
Set-TextValue $ws.Range('E50') '  +1.23%  '
Set-TextValue $ws.Range('D51') '0.0270'
Set-TextValue $ws.Range('E51') '  +3.15%  '
